# "Added some pics and updated spreadsheet"
# Updates the CBLeaves notes sheet: extends the F8 note, shortens the
# "Larimar/amber tree?" entry into separate Larimar/Amber rows, and fills
# in season-count labels for the VMV (F/G) and (shifted) SVE (I/J) mini
# tables, plus a brand-new RSV mini table in columns L/M.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Expand the cherry-blossom note.
$ws.Range("F8").Value = "Those big pink trees randomly around; cherry blossoms? 1 week per Spring"

# VMV list (column F) gains season-count labels in column G.
# SVE list moves from column J into column I, one column left.
# Column J is repurposed to hold season-count labels.
# A brand-new RSV list (name + count) lands in columns L/M.
$ws.Range("G17").Value = 24
$ws.Range("I17").Value = "SVE"
$ws.Range("J17").Value = 20
$ws.Range("L17").Value = "RSV"
$ws.Range("M17").Value = 23
$ws.Range("N17").ClearContents()

$ws.Range("G18").Value = "4 seasons"
$ws.Range("I18").Value = "Douglas fir tree"
$ws.Range("J18").Value = "4 seasons"
$ws.Range("L18").Value = "Cherry pluot"
$ws.Range("M18").Value = 3

$ws.Range("G19").Value = "4 seasons"
$ws.Range("I19").Value = "Persimmon tree"
$ws.Range("J19").Value = "3 seasons"
$ws.Range("L19").Value = "Mtn plumcot"
$ws.Range("M19").Value = 3

$ws.Range("G20").Value = "3 seasons"
$ws.Range("I20").Value = "Pear tree"
$ws.Range("J20").Value = "3 seasons"
$ws.Range("L20").Value = "Desert tangelo"
$ws.Range("M20").Value = 3

$ws.Range("G21").Value = "3 seasons"
$ws.Range("I21").Value = "Nectarine tree"
$ws.Range("J21").Value = "3 seasons"
$ws.Range("L21").Value = "Paradise rangpur"
$ws.Range("M21").Value = 2

# "Larimar/amber tree?" splits into two rows: Larimar and Amber.
$ws.Range("F22").Value = "Larimar"
$ws.Range("G22").Value = "1 season"
$ws.Range("I22").Value = "Money tree?"
$ws.Range("J22").Value = "4 seasons"
$ws.Range("L22").Value = "Tropi ugli"
$ws.Range("M22").Value = 3

$ws.Range("F23").Value = "Amber"
$ws.Range("G23").Value = "1 season"
$ws.Range("I23").Value = "Birch tree"
$ws.Range("J23").Value = "3 seasons"
$ws.Range("L23").Value = "Ember blood lime"
$ws.Range("M23").Value = 3

$ws.Range("F24").Value = "Casolatier"
$ws.Range("G24").Value = "4 seasons"
$ws.Range("L24").Value = "Highland jostaberry"
$ws.Range("M24").Value = 3

$ws.Range("F25").Value = "Trellis grape"
$ws.Range("G25").Value = "1 season"
$ws.Range("L25").Value = "North limequat"
$ws.Range("M25").Value = 3

$ws.Range("F26").Value = "Nevaril bush"
$ws.Range("G26").Value = "1 season"

$ws.Range("F27").Value = "Rosemary bush"
$ws.Range("G27").Value = "1 season"

$ws.Range("F28").Value = "Sea buckthorn bush"
$ws.Range("G28").Value = "1 season"

# Match the author's final selection.
$ws.Range("J30").Select()
